# Move the "_GoBack" bookmark.
#
# In the original document the (hidden) "_GoBack" bookmark sits in the empty
# paragraph that follows the "Final Index Creation" heading, right near the
# end of the document. The edit relocates it to the end of the "Evaluation"
# paragraph that closes out the Logistic Regression section (i.e. right
# before the blank paragraph that precedes the "Random Forest Classifier"
# heading) - no visible text changes, just the bookmark's position.

$d = $word.ActiveDocument

# --- Find the destination paragraph -----------------------------------
# It is identified by its text being "Evaluation" immediately followed by
# a blank paragraph and then the "Random Forest Classifier" heading.
$needle = "Evaluation" + [char]13 + [char]13 + "Random Forest Classifier"
$span = $d.Content
$found = $span.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $destPos = $span.Start + ("Evaluation").Length
    $dest = $d.Range($destPos, $destPos)

    # Word's bookmark engine here can't reliably anchor a bookmark to a
    # *collapsed* range that lands exactly on a paragraph-mark boundary next
    # to an empty paragraph. Work around it by temporarily widening the
    # destination with a throw-away character, bookmarking that, and then
    # deleting the character again (the bookmark collapses back down with
    # it, the way Word bookmarks always shrink to fit a deletion).
    $dest.InsertAfter([char]1)

    # Re-adding a bookmark with the same name removes/replaces any existing
    # bookmark of that name elsewhere in the document, so this both creates
    # the bookmark here and removes it from its old location in one step.
    $d.Bookmarks.Add("_GoBack", $dest) | Out-Null

    $shrink = $d.Range($destPos, $destPos + 1)
    $shrink.Text = ""
}
